$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 208
$ws.Range("B208").Value = 7404216
$ws.Range("F208").Value = 'Independiente Santa Fe'
$ws.Range("G208").Value = 'Once Caldas'
$ws.Range("H208").Value = 0
$ws.Range("J208").Value = 'A'
$ws.Range("K208").Value = 1.85
$ws.Range("L208").Value = 3.1
$ws.Range("M208").Value = 4.2
$ws.Range("N208").Value = 2.25
$ws.Range("P208").Value = 3.3
$ws.Range("Q208").Value = -0.25
$ws.Range("R208").Value = 1.9
$ws.Range("S208").Value = 1.9
$ws.Range("U208").Value = 1.925
$ws.Range("V208").Value = 1.925
$ws.Range("X208").Value = -1
$ws.Range("Y208").Value = 2.3
$ws.Range("Z208").Value = -1
$ws.Range("AA208").Value = 0.8999999999999999
$ws.Range("AC208").Value = 0.925

# Row 209
$ws.Range("B209").Value = 7404217
$ws.Range("F209").Value = 'Alianza Petrolera'
$ws.Range("G209").Value = 'Deportivo Pereira'
$ws.Range("H209").Value = 2
$ws.Range("J209").Value = 'H'
$ws.Range("K209").Value = 1.95
$ws.Range("L209").Value = 3.2
$ws.Range("M209").Value = 3.75
$ws.Range("N209").Value = 1.95
$ws.Range("O209").Value = 3.2
$ws.Range("P209").Value = 4.75
$ws.Range("Q209").Value = -0.5
$ws.Range("R209").Value = 1.925
$ws.Range("S209").Value = 1.875
$ws.Range("T209").Value = 2
$ws.Range("U209").Value = 1.825
$ws.Range("W209").Value = 0.95
$ws.Range("X209").Value = -1
$ws.Range("Z209").Value = 0.925
$ws.Range("AA209").Value = -1
$ws.Range("AB209").Value = 0.825
$ws.Range("AC209").Value = -1

# Row 211
$ws.Range("B211").Value = 7404212
$ws.Range("F211").Value = 'Envigado FC'
$ws.Range("G211").Value = 'Deportivo Pasto'
$ws.Range("H211").Value = 1
$ws.Range("J211").Value = 'D'
$ws.Range("K211").Value = 2.6
$ws.Range("L211").Value = 2.875
$ws.Range("M211").Value = 2.8
$ws.Range("N211").Value = 2.8
$ws.Range("P211").Value = 2.625
$ws.Range("Q211").Value = 0
$ws.Range("R211").Value = 1.975
$ws.Range("T211").Value = 2.5
$ws.Range("U211").Value = 2.025
$ws.Range("V211").Value = 1.825
$ws.Range("W211").Value = -1
$ws.Range("X211").Value = 2.2
$ws.Range("Z211").Value = 0
$ws.Range("AA211").Value = -0
$ws.Range("AB211").Value = -1
$ws.Range("AC211").Value = 0.825

# Row 212
$ws.Range("B212").Value = 7404214
$ws.Range("F212").Value = 'Boyaca Chico'
$ws.Range("G212").Value = 'Deportivo Cali'
$ws.Range("H212").Value = 1
$ws.Range("J212").Value = 'D'
$ws.Range("K212").Value = 3.2
$ws.Range("M212").Value = 2.2
$ws.Range("N212").Value = 3.6
$ws.Range("O212").Value = 3
$ws.Range("P212").Value = 2.25
$ws.Range("Q212").Value = 0.25
$ws.Range("R212").Value = 1.95
$ws.Range("T212").Value = 2.25
$ws.Range("U212").Value = 1.875
$ws.Range("V212").Value = 1.975
$ws.Range("X212").Value = 2
$ws.Range("Y212").Value = -1
$ws.Range("Z212").Value = 0.475
$ws.Range("AA212").Value = -0.5
$ws.Range("AB212").Value = -0.5
$ws.Range("AC212").Value = 0.4875

# Row 240
$ws.Range("B240").Value = 7528135
$ws.Range("F240").Value = 'Independiente Medellin'
$ws.Range("G240").Value = 'America de Cali'
$ws.Range("H240").Value = 2
$ws.Range("I240").Value = 1
$ws.Range("K240").Value = 2.15
$ws.Range("M240").Value = 3.4
$ws.Range("N240").Value = 2.375
$ws.Range("O240").Value = 3.3
$ws.Range("P240").Value = 3.1
$ws.Range("Q240").Value = -0.25
$ws.Range("R240").Value = 2
$ws.Range("S240").Value = 1.8
$ws.Range("U240").Value = 1.975
$ws.Range("V240").Value = 1.825
$ws.Range("W240").Value = 1.375
$ws.Range("Z240").Value = 1
$ws.Range("AB240").Value = 0.9750000000000001

# Row 241
$ws.Range("B241").Value = 7528603
$ws.Range("F241").Value = 'Junior'
$ws.Range("G241").Value = 'Deportes Tolima'
$ws.Range("H241").Value = 4
$ws.Range("I241").Value = 2
$ws.Range("K241").Value = 1.95
$ws.Range("M241").Value = 4
$ws.Range("N241").Value = 1.909
$ws.Range("O241").Value = 3.75
$ws.Range("P241").Value = 3.8
$ws.Range("Q241").Value = -0.5
$ws.Range("R241").Value = 1.9
$ws.Range("S241").Value = 1.9
$ws.Range("U241").Value = 1.85
$ws.Range("V241").Value = 1.95
$ws.Range("W241").Value = 0.909
$ws.Range("Z241").Value = 0.8999999999999999
$ws.Range("AB241").Value = 0.8500000000000001

# Row 373
$ws.Range("B373").Value = 7658955
$ws.Range("E373").Value = 45381.84722222222
$ws.Range("F373").Value = 'Boyaca Chico'
$ws.Range("G373").Value = 'Jaguares de Cordoba'
$ws.Range("K373").Value = 1.909
$ws.Range("L373").Value = 3.4
$ws.Range("M373").Value = 4
$ws.Range("N373").Value = 1.6
$ws.Range("O373").Value = 4
$ws.Range("P373").Value = 5.75
$ws.Range("Q373").Value = -0.75
$ws.Range("R373").Value = 1.775
$ws.Range("S373").Value = 2.1
$ws.Range("U373").Value = 1.95
$ws.Range("V373").Value = 1.9

# Row 374
$ws.Range("B374").Value = 7658952
$ws.Range("E374").Value = 45381.9375
$ws.Range("F374").Value = 'Fortaleza'
$ws.Range("G374").Value = 'Millonarios'
$ws.Range("K374").Value = 2.6
$ws.Range("L374").Value = 3.25
$ws.Range("M374").Value = 2.75
$ws.Range("N374").Value = 2.05
$ws.Range("O374").Value = 3.3
$ws.Range("P374").Value = 3.75
$ws.Range("Q374").Value = -0.25
$ws.Range("R374").Value = 1.775
$ws.Range("S374").Value = 2.1
$ws.Range("T374").Value = 2
$ws.Range("U374").Value = 1.85
$ws.Range("V374").Value = 2

# Row 375
$ws.Range("B375").Value = 7658950
$ws.Range("E375").Value = 45382.66666666666
$ws.Range("F375").Value = 'Independiente Santa Fe'
$ws.Range("G375").Value = 'Patriotas FC'
$ws.Range("K375").Value = 1.5
$ws.Range("L375").Value = 4
$ws.Range("M375").Value = 6.5
$ws.Range("N375").Value = 1.6
$ws.Range("O375").Value = 3.8
$ws.Range("P375").Value = 6
$ws.Range("Q375").Value = -0.75
$ws.Range("R375").Value = 1.775
$ws.Range("S375").Value = 2.1

# Row 376
$ws.Range("B376").Value = 7658951
$ws.Range("E376").Value = 45382.75694444445
$ws.Range("F376").Value = 'Deportivo Pereira'
$ws.Range("G376").Value = 'Atletico Bucaramanga'
$ws.Range("K376").Value = 1.833
$ws.Range("L376").Value = 3.25
$ws.Range("M376").Value = 4.5
$ws.Range("N376").Value = 1.727
$ws.Range("O376").Value = 3.4
$ws.Range("P376").Value = 5.75
$ws.Range("Q376").Value = -0.75
$ws.Range("R376").Value = 1.975
$ws.Range("S376").Value = 1.875
$ws.Range("U376").Value = 1.8
$ws.Range("V376").Value = 2.05

# Row 377
$ws.Range("B377").Value = 7658949
$ws.Range("E377").Value = 45382.84722222222
$ws.Range("F377").Value = 'Independiente Medellin'
$ws.Range("G377").Value = 'America de Cali'
$ws.Range("K377").Value = 2.4
$ws.Range("L377").Value = 3.25
$ws.Range("M377").Value = 2.9
$ws.Range("N377").Value = 2.45
$ws.Range("O377").Value = 3.25
$ws.Range("P377").Value = 3
$ws.Range("Q377").Value = -0.25
$ws.Range("R377").Value = 2.1
$ws.Range("S377").Value = 1.775
$ws.Range("T377").Value = 2.25
$ws.Range("U377").Value = 1.9
$ws.Range("V377").Value = 1.95

# Row 378
$ws.Range("B378").Value = 7658954
$ws.Range("E378").Value = 45382.9375
$ws.Range("F378").Value = 'Deportivo Cali'
$ws.Range("G378").Value = 'Aguilas Doradas'
$ws.Range("K378").Value = 2.1
$ws.Range("M378").Value = 3.5
$ws.Range("N378").Value = 2.15
$ws.Range("O378").Value = 3.25
$ws.Range("P378").Value = 3.6
$ws.Range("R378").Value = 1.825
$ws.Range("S378").Value = 2.025
$ws.Range("T378").Value = 2.25
$ws.Range("U378").Value = 2.1
$ws.Range("V378").Value = 1.775

# Delete trailing rows 379-382 (matches removed in this update)
$ws.Range("A379:AC382").Delete()
